$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.472.89"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "1.569.91"
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.63"
$ws.Range("E5").Value = "  +1.36%  "
$ws.Range("E6").Value = "  -1.18%  "
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.19"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0592"
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").Value = "1.793.70"
$ws.Range("E12").Value = "  -0.89%  "
$ws.Range("D13").Value = "1.568.30"
$ws.Range("E13").Value = "  -1.01%  "
$ws.Range("E14").Value = "  -1.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.518"
$ws.Range("E15").Value = "  -2.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.79"
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("D17").Value = "27.477.49"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "214.30"
$ws.Range("E18").Value = "  -2.15%  "
$ws.Range("D19").Value = "0.0₃0691"
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("E22").Value = "  -0.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.56"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("E24").Value = "  +1.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.48"
$ws.Range("E25").Value = "  -1.76%  "
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("E28").Value = "  -0.52%  "
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("E30").Value = "  +0.43%  "
$ws.Range("E31").Value = "  +0.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.20"
$ws.Range("D33").Value = "1.381.18"
$ws.Range("E34").Value = "  +1.74%  "
$ws.Range("E35").Value = "  +1.00%  "
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.953"
$ws.Range("E36").Value = "  -2.60%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.31"
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("E38").Value = "  +0.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.543"
$ws.Range("E39").Value = "  +0.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.828"
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.980"
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("E43").Value = "  +1.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.22"
$ws.Range("E44").Value = "  +1.23%  "
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("E46").Value = "  +0.63%  "
$ws.Range("D47").Value = "1.705.50"
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.50"
$ws.Range("E48").Value = "  -3.18%  "
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0497"
$ws.Range("E50").Value = "  -0.54%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0957"
$ws.Range("E51").Value = "  -1.72%  "
